$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) and E (Volume 1h) to Text format so that
# numeric-looking strings (e.g. "347.71") are stored as text, matching
# the original inline-string cell contents rather than being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.868.95"
$ws.Range("D3").Value = "2.115.04"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "347.71"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "0.5184"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "0.4461"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("E9").Value = "  +3.54%  "
$ws.Range("D10").Value = "0.09354"
$ws.Range("E10").Value = "  +4.23%  "
$ws.Range("D11").Value = "1.180"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "25.13"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").Value = "2.106.62"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "8.389"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "102.60"
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("D17").Value = "0.00001165"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "21.56"
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").Value = "0.06672"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "6.299"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "29.911.59"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "12.70"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "2.326"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "2.357.54"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").Value = "22.12"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").Value = "2.556"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").Value = "162.54"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "134.02"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "1.153"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").Value = "1.793"
$ws.Range("E32").Value = "  +9.16%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "6.244"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "3.979"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "6.466"
$ws.Range("E36").Value = "  +5.51%  "
$ws.Range("D37").Value = "10.89"
$ws.Range("E37").Value = "  +7.62%  "
$ws.Range("D38").Value = "0.02600"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").Value = "12.68"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "0.7016"
$ws.Range("E41").Value = "  +3.20%  "
$ws.Range("D42").Value = "1.349"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "0.2240"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "0.6848"
$ws.Range("E44").Value = "  +6.62%  "
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("D46").Value = "2.363"
$ws.Range("E46").Value = "  +3.66%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "0.00000000360"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "3.642"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "1.233"
$ws.Range("E50").Value = "  +5.87%  "
$ws.Range("D51").Value = "1.225"
$ws.Range("E51").Value = "  +0.65%  "

# Restore default cell style (removes the temporary Text number format)
$ws.Range("D2:E51").Style = "Normal"

